# Apply "Added support for priors" changes to the experiment template workbook.

$wb = $excel.ActiveWorkbook

$wsExpDesc  = $wb.Worksheets.Item("experiment_description")
$wsExpSpec  = $wb.Worksheets.Item("experiment_specification")
$wsRunDesc  = $wb.Worksheets.Item("run_description")
$wsRunSpec  = $wb.Worksheets.Item("run_specification")
$wsHeurDesc = $wb.Worksheets.Item("heuristics_description")

# --- experiment_description: add prior_transition / prior_los columns ---
$wsExpDesc.Range("E1").Value = "prior_transition"
$wsExpDesc.Range("F1").Value = "prior_los"

# All existing experiments (rows 2-10 and 12-14) get "none" for both new columns
$wsExpDesc.Range("E2:F10").Value = "none"
$wsExpDesc.Range("E12:F14").Value = "none"

# Row 11 (experiment_id 10) becomes the new "Ferguson Wuhan model" experiment
$wsExpDesc.Range("B11").Value = "Ferguson Wuhan model"
$wsExpDesc.Range("C11").Value = "ferguson"
$wsExpDesc.Range("D11").Value = "4;6"
$wsExpDesc.Range("E11").Value = "wuhan"
$wsExpDesc.Range("F11").Value = "wuhan"

# --- experiment_specification: experiment 10 rows get new splitting values ---
$wsExpSpec.Range("C29").Value = "length_of_stay_simple_two_weeks"
$wsExpSpec.Range("D29").Value = "age_official"
$wsExpSpec.Range("C30").Value = "length_of_stay_simple_two_weeks"
$wsExpSpec.Range("D30").Value = "age_official"

# --- sheet view / selection adjustments to mirror the saved workbook state ---
[void]$wsExpSpec.Activate()
try { $excel.ActiveWindow.ScrollRow = 6 } catch { }
[void]$wsExpSpec.Range("C31").Select()

[void]$wsRunDesc.Activate()
[void]$wsRunDesc.Range("B9").Select()

[void]$wsRunSpec.Activate()
[void]$wsRunSpec.Range("B7").Select()

[void]$wsHeurDesc.Activate()
[void]$wsHeurDesc.Range("B7").Select()

[void]$wsExpDesc.Activate()
[void]$wsExpDesc.Range("F11").Select()

$wb.Save()
